$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $val) {
    $r = $ws.Range($addr)
    $r.NumberFormat = "@"
    $r.Value = $val
    $r.ClearFormats()
}

$ws.Range("D2").Value = "27.723.10"
$ws.Range("E2").Value = "  +0.05%  "
$ws.Range("D3").Value = "1.904.75"
$ws.Range("E3").Value = "  +0.53%  "
Set-TextValue "D4" "0.9992"
$ws.Range("E4").Value = "  -0.18%  "
Set-TextValue "D5" "312.69"
$ws.Range("E5").Value = "  -0.12%  "
$ws.Range("E6").Value = "  -0.11%  "
Set-TextValue "D7" "0.5210"
$ws.Range("E7").Value = "  +7.25%  "
Set-TextValue "D8" "0.3778"
$ws.Range("E8").Value = "  -0.36%  "
$ws.Range("E9").Value = "  -1.31%  "
Set-TextValue "D10" "21.28"
$ws.Range("E10").Value = "  +3.55%  "
Set-TextValue "D12" "0.07638"
$ws.Range("E12").Value = "  -0.71%  "
$ws.Range("D13").Value = "1.878.46"
$ws.Range("E13").Value = "  -0.64%  "
Set-TextValue "D14" "5.448"
$ws.Range("E14").Value = "  -0.43%  "
Set-TextValue "D15" "92.10"
$ws.Range("E15").Value = "  +1.28%  "
Set-TextValue "D16" "0.9997"
$ws.Range("E16").Value = "  -0.20%  "
Set-TextValue "D17" "0.000008705"
$ws.Range("E17").Value = "  -1.17%  "
Set-TextValue "D18" "0.9999"
$ws.Range("E18").Value = "  -0.06%  "
$ws.Range("D19").Value = "27.757.97"
$ws.Range("E19").Value = "  +0.02%  "
Set-TextValue "D20" "14.50"
$ws.Range("E20").Value = "  +0.20%  "
Set-TextValue "D21" "5.143"
$ws.Range("E21").Value = "  +0.47%  "
$ws.Range("D22").Value = "2.140.50"
$ws.Range("E22").Value = "  +1.78%  "
Set-TextValue "D23" "10.83"
$ws.Range("E23").Value = "  +0.79%  "
Set-TextValue "D24" "6.583"
$ws.Range("E24").Value = "  -0.35%  "
Set-TextValue "D25" "153.36"
$ws.Range("E25").Value = "  -0.21%  "
Set-TextValue "D26" "1.879"
$ws.Range("E26").Value = "  -1.37%  "
$ws.Range("B27").Value = "LidoDAOToken"
$ws.Range("C27").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
Set-TextValue "D27" "2.168"
$ws.Range("E27").Value = "  +0.96%  "
$ws.Range("B28").Value = "EthereumClassic"
$ws.Range("C28").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
Set-TextValue "D28" "18.32"
$ws.Range("E28").Value = "  -0.18%  "
$ws.Range("E29").Value = "  -1.25%  "
Set-TextValue "D30" "4.864"
$ws.Range("E30").Value = "  -0.75%  "
Set-TextValue "D31" "0.08983"
$ws.Range("E31").Value = "  +0.83%  "
Set-TextValue "D32" "4.855"
$ws.Range("E32").Value = "  +4.52%  "
Set-TextValue "D33" "3.176"
$ws.Range("E33").Value = "  +0.70%  "
Set-TextValue "D34" "1.230"
$ws.Range("E34").Value = "  +0.50%  "
Set-TextValue "D35" "0.7732"
$ws.Range("E35").Value = "  +1.14%  "
Set-TextValue "D36" "2.635"
$ws.Range("E36").Value = "  +4.47%  "
Set-TextValue "D37" "0.02084"
$ws.Range("E37").Value = "  +2.10%  "
Set-TextValue "D38" "3.068"
$ws.Range("E38").Value = "  +2.71%  "
$ws.Range("E39").Value = "  +0.14%  "
Set-TextValue "D40" "0.5516"
$ws.Range("E40").Value = "  +0.95%  "
Set-TextValue "D41" "0.05291"
$ws.Range("E41").Value = "  +0.32%  "
Set-TextValue "D42" "6.675"
$ws.Range("E42").Value = "  -3.54%  "
Set-TextValue "D43" "114.60"
$ws.Range("E43").Value = "  +4.17%  "
$ws.Range("E44").Value = "  +0.94%  "
Set-TextValue "D45" "0.1510"
$ws.Range("E45").Value = "  -0.45%  "
Set-TextValue "D46" "0.4808"
$ws.Range("E46").Value = "  +0.53%  "
Set-TextValue "D47" "10.38"
$ws.Range("E47").Value = "  -2.13%  "
Set-TextValue "D48" "0.9994"
$ws.Range("E48").Value = "  -0.12%  "
Set-TextValue "D49" "1.615"
$ws.Range("E49").Value = "  -1.32%  "
Set-TextValue "D50" "66.78"
$ws.Range("E50").Value = "  -0.80%  "
Set-TextValue "D51" "0.05986"
$ws.Range("E51").Value = "  -1.12%  "

Write-Output "applied 97 cell changes"